$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 0.0004861111111111111
$ws.Range("K2").Value = 3474
$ws.Range("L2").Value = 0.006948
